# Recreate the four PowerPoint "sections" (Slide Sorter groupings) that the
# author set up on this deck:
#   1. Empty Section 1          -> no slides
#   2. Section with Slides 1    -> slides 1-4 (sldId 256-259)
#   3. Section with Slides 2    -> slides 5-8 (sldId 260-263)
#   4. Empty Section 2          -> no slides
#
# PowerPoint's real COM model exposes this through
# Presentation.SectionProperties, with:
#   AddSection(index, name)       -> inserts a (normally empty) section at
#                                     the given 1-based position
#   AddBeforeSlide(slideIdx, name)-> the "Add Section" command you get from
#                                     right-clicking a slide in Slide Sorter;
#                                     it starts a new, named section right
#                                     before that slide and absorbs the
#                                     slides that would otherwise trail the
#                                     previous section.

$p = $ppt.ActivePresentation
$sections = $p.SectionProperties

# First, an empty section right at the very top of the deck.
[void]$sections.AddSection(1, "Empty Section 1")

# Then split the eight slides into two four-slide sections.
[void]$sections.AddBeforeSlide(1, "Section with Slides 1")
[void]$sections.AddBeforeSlide(5, "Section with Slides 2")

# Finally, an empty section after the last slide.
[void]$sections.AddSection($sections.Count + 1, "Empty Section 2")
